$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.187.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.53%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.835.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.35%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.9989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'242.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.43%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.6593"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.11%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.07425"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.58%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.2928"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.51%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'22.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.32%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("B11").Value = "'WrappedEther"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'1.987.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +8.41%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").Value = "'TRON"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.07760"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.61%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'4.987"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.23%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.6656"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.76%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'82.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.82%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "'Uniswap"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'6.128"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.27%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = "'ShibaInu"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.000008595"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +5.07%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = "'WrappedBTC"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'29.182.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.52%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("B19").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'2.063.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.04%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'226.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.28%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'12.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.43%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  -0.06%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'7.121"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.76%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.9999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.04%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'159.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.77%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'8.606"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.35%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.1398"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.38%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'17.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.11%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  +1.06%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'4.117"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.49%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'4.047"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.31%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'1.194"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.32%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.05258"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.40%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.867"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.35%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.7354"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.49%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'1.145"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.79%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'2.656"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.05%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'1.302.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.60%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.01795"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.57%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'2.735"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.31%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.9219"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.11%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.08818"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +14.80%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'6.034"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.39%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = "'PaxDollar"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.9994"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.01%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'Quant"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'102.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.94%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'RocketPoolETH"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1.956.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.16%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'Mantle"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.5141"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.68%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "'Aave"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'63.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.77%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'RenderToken"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'1.753"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.51%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "'BabyDogeCoin"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.00000000120"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.66%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  -1.10%  "
$ws.Range("E51").Style = "Normal"
